$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.535.22'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '1.642.01'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.22'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3763'
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.68'
$ws.Range("E8").Value = '  +3.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3669'
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.281'
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08198'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9995'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.13'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.681'
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001284'
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.467'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("D17").Value = '1.644.82'
$ws.Range("E17").Value = '  +3.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.26'
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06936'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.35'
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.587'
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9976'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").Value = '23.539.92'
$ws.Range("E23").Value = '  +1.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.073'
$ws.Range("E25").Value = '  +4.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.419'
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.35'
$ws.Range("E27").Value = '  +1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.31'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.342'
$ws.Range("E29").Value = '  +2.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.44'
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.373'
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '1.828.65'
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.851'
$ws.Range("E33").Value = '  +1.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9763'
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02852'
$ws.Range("E35").Value = '  +6.56%  '
$ws.Range("E36").Value = '  +2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.07429'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2557'
$ws.Range("E38").Value = '  +2.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.205'
$ws.Range("E39").Value = '  +1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08917'
$ws.Range("E40").Value = '  +1.46%  '
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7141'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.60'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.33'
$ws.Range("E44").Value = '  +7.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6583'
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.357'
$ws.Range("E46").Value = '  +3.17%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.041'
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9984'
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08011'
$ws.Range("E49").Value = '  +1.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '130.03'
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("E51").Value = '  +1.60%  '
